$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2. This shifts the existing
# row 2 (and everything below) down by one, turning old row 2 into row 3.
$ws.Rows.Item(2).Insert()

# Make sure the date-like text values we are about to write are kept as
# plain text instead of being auto-converted into Excel date serials.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("AA2").NumberFormat = "@"

# Populate the new row 2 with the new observation record.
$ws.Range("A2").Value = 131092986
$ws.Range("B2").Value = 97878
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 221945
$ws.Range("F2").Value = "Revlummer"
$ws.Range("G2").Value = "Lycopodium annotinum"
$ws.Range("H2").Value = "L."
$ws.Range("P2").Value = "254, Upl"
$ws.Range("Q2").Value = 693755
$ws.Range("R2").Value = 6663402
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = "Stockholm"
$ws.Range("U2").Value = "Norrtälje"
$ws.Range("V2").Value = "Uppland"
$ws.Range("W2").Value = "Edebo"
$ws.Range("Y2").Value = "2024-09-24"
$ws.Range("AA2").Value = "2024-09-24"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = "Anton Ringbom"
$ws.Range("AX2").Value = "Anton Ringbom, Alexandra Östberg"
